$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the style used by the other
# header cells (bold font, thin border, centered/top alignment = G1's style).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for the Save column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
